$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date on Francisco Medeiros' row (D2): was 80832, should be 44308
$ws.Range("D2").Value = 44308

# Add a new row (row 3) for Jefferson André, same Curso/Carga as row 2
$ws.Range("A3").Value = "Jefferson André"
$ws.Range("B3").Value = "Example"
$ws.Range("C3").Value = "5 hours"

# Copy the date's number formatting from D2, then set the new date value
$ws.Range("D2").Copy($ws.Range("D3"))
$ws.Range("D3").Value = 44308

# Set the e-mail text, then turn it into a mailto hyperlink (matching E2's
# existing hyperlink), and re-apply the Hyperlink cell style
$ws.Range("E3").Value = "jefferson.andre96@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:jefferson.andre96@gmail.com")
$ws.Range("E3").Style = $ws.Range("E2").Style

# Update the active selection to D9
$ws.Range("D9").Select() | Out-Null
